# Apply the Spotify-clone spreadsheet fixes:
#  1) Fix text typos / normalize comma+space separators inside the quoted
#     song-list / artist-list cells (and the "familiar" spelling fix).
#  2) Remove the blank spacer column A and the blank spacer rows so the
#     sheet is tight around the two data tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix cell text (typos / spacing) -----------------------------------
$ws.Range("G4").Value  = '"Soul For Us", "Magic Circus", "Diamond Power", "Thang Of Thunder"'
$ws.Range("E5").Value  = "familiar"
$ws.Range("G5").Value  = '"Home Forever", "Words Of Her Life", "Reflections Of Magic", "Honey, Let''s Be Silly"'
$ws.Range("G7").Value  = '"Dance With Her Own", "Without My Streets", "Celebration Of More"'
$ws.Range("E12").Value = '"Soul For Us", "Reflections Of Magic", "Dance With Her Own"'
$ws.Range("E13").Value = '"Troubles Of My Inner Fire", "Time Fireworks"'
$ws.Range("E14").Value = '"Magic Circus", "Honey, So Do I", "Sweetie, Let''s Go Wild", "She Knows"'
$ws.Range("E15").Value = '"Fantasy For Me", "Celebration Of More", "Rock His Everything", "Home Forever", "Diamond Power", "Honey, Let''s Be Silly"'
$ws.Range("E16").Value = '"Thang Of Thunder", "Words Of Her Life", "Without My Streets"'

# --- 2) Remove blank spacer rows/column -----------------------------------
# Remove the trailing blank rows under the second table first (rows 17-21),
# then the blank rows between the two tables (rows 8-9), then the blank
# rows above the first table (rows 1-2) -- bottom to top so row numbers of
# not-yet-deleted rows stay stable.
$ws.Rows("17:21").Delete()
$ws.Rows("8:9").Delete()
$ws.Rows("1:2").Delete()

# Remove the narrow blank spacer column A.
$ws.Columns("A").Delete()
